# Update the test credentials on Sheet1.
# Row 1 (A1:B1) stays as the "username"/"password" header.
# Rows 2-6 are rewritten to use the new manager login/password pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUser = "mngr429183"
$newPass = "Avamyze"

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = $newUser
    $ws.Cells.Item($r, 2).Value = $newPass
}

# Matches the resulting selection left behind in the saved file.
$ws.Range("A8").Select()
